$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.966.07"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.819.05"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'310.97"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.4675"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "'0.3663"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "'0.07348"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'0.8732"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "'20.28"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "1.802.31"
$ws.Range("E12").Value = "  -3.04%  "
$ws.Range("D13").Value = "'5.413"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'0.07157"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "'91.51"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "'0.000008738"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "'14.66"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "26.995.60"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'5.290"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "2.056.73"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").Value = "'151.09"
$ws.Range("D27").Value = "'18.31"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "'2.140"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").Value = "'5.237"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").Value = "'116.95"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "'0.08889"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'0.7555"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "'1.158"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'4.500"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'2.944"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").Value = "'1.096"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'0.05304"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").Value = "'2.969"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").Value = "'2.382"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "'7.176"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "'0.5294"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "'0.1653"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'8.455"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'0.4887"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "'10.47"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "'1.664"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "'103.12"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "'0.06296"
$ws.Range("E51").Value = "  +0.11%  "
